$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the sheet that no longer exists (TEST_STANDARD_ROW) and rename
#    TEST_SCALAR_INPUT -> TEST_SCALAR (tracked rename), then reorder so the
#    renamed sheet sits right after TXL_SHEET_TRACKER (where TEST_SHEET used
#    to be), matching the new tab order.
# ---------------------------------------------------------------------------
$wsStandardRow = $wb.Worksheets.Item("TEST_STANDARD_ROW")
$wsStandardRow.Delete()

$wsScalarInput = $wb.Worksheets.Item("TEST_SCALAR_INPUT")
$wsScalarInput.Name = "TEST_SCALAR"

$wsTestSheet = $wb.Worksheets.Item("TEST_SHEET")
$wsScalar = $wb.Worksheets.Item("TEST_SCALAR")
$wsScalar.Move($wsTestSheet)

# ---------------------------------------------------------------------------
# 2. Rebuild the defined names so they reference the renamed sheet and match
#    the tracker's new naming scheme.
# ---------------------------------------------------------------------------
$existingNames = @()
foreach ($n in $wb.Names) {
    $existingNames += $n.Name
}
foreach ($nm in $existingNames) {
    $wb.Names.Item($nm).Delete()
}

$wb.Names.Add("TEST_SCALAR__dollar", "=TEST_SCALAR!`$B`$2")
$wb.Names.Add("TEST_SCALAR__kwh", "=TEST_SCALAR!`$B`$3")
$wb.Names.Add("TEST_SCALAR__value", "=TEST_SCALAR!`$B:`$B")
$wb.Names.Add("TEST_SCALAR__var_name", "=TEST_SCALAR!`$A:`$A")

# ---------------------------------------------------------------------------
# 3. Update the TXL_SHEET_TRACKER sheet: the tracked row now points at
#    TEST_SCALAR, the live SHEET() lookup formula (which referenced the
#    now-deleted TEST_STANDARD_ROW sheet) is replaced with its resolved,
#    static sheet-index value, and the first column is narrowed.
# ---------------------------------------------------------------------------
$wsTracker = $wb.Worksheets.Item("TXL_SHEET_TRACKER")
$wsTracker.Range("A2").Value = "TEST_SCALAR"
$wsTracker.Range("B2").Value = "Sheet to test standard row sheet"
$wsTracker.Range("D2").Value = 2

$wsTracker.Columns.Item(1).ColumnWidth = 13.285714285714286

# ---------------------------------------------------------------------------
# 4. Restore per-sheet selections to match the new layout.
# ---------------------------------------------------------------------------
$wsTracker.Activate()
$wsTracker.Range("B32").Select()

$wsScalar = $wb.Worksheets.Item("TEST_SCALAR")
$wsScalar.Activate()
$wsScalar.Range("D32").Select()

$wsTestSheet = $wb.Worksheets.Item("TEST_SHEET")
$wsTestSheet.Activate()
$wsTestSheet.Range("E35").Select()

$wsTracker.Activate()
